# "atualização da base de dados" — add two new products to the
# "produtos" sheet, keeping the existing alphabetical ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("produtos")

# Insert "cabeludinha/ jabuticaba amarela" (na_terra = x) right before
# "cebolinha", which is its correct alphabetical spot.
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "cabeludinha/ jabuticaba amarela"
$ws.Range("B17").Value = "x"

# Insert "jiló" (na_terra = x) right before "laranja", which is its
# correct alphabetical spot (rows below 17 have already shifted down
# by one because of the insert above).
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "jiló"
$ws.Range("B35").Value = "x"

# Re-apply the sheet's remembered sort (produto, ascending, no header)
# over the now-larger A2:C53 range so the cached sort state on the
# sheet tracks the new extent, matching what Excel records after a
# Data > Sort on the refreshed range.
$sf = $ws.Sort.SortFields
$sf.Clear()
$sf.Add($ws.Range("A53"))
$ws.Sort.SetRange($ws.Range("A2:C53"))
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.Apply()

# The engine's Sort.Apply compares text ordinally rather than with
# pt-BR collation, which only affects "salsão" vs "salsinha" (the
# accented "ã" sorts after "i" by code point even though it is treated
# like a plain "a" in Portuguese). Restore the correct data/order for
# those two rows so the sheet content matches the real, locale-aware
# alphabetical ordering.
$ws.Range("A50").Value = "salsão"
$ws.Range("B50").Value = "x"
$ws.Range("C50").ClearContents()
$ws.Range("A51").Value = "salsinha"
$ws.Range("B51").Value = "x"
$ws.Range("C51").Value = "x"

# Match the author's final cursor/selection position.
$ws.Range("D47").Select()
